$wb = $excel.ActiveWorkbook

# Sheet 1: GNG_TO...
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504778363739493"
$ws1.Range("B2").Value = "go_stims-1650477836339949.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778363559551.csv"
$ws1.Range("B4").Value = "go_stims-16504778363579533.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778363719842.csv"

# Sheet 2: NB_TO...
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778385419507"
$ws2.Range("B2").Value = "OB-16504778368229847.csv"
$ws2.Range("B3").Value = "ZB-match_0-16504778366379874.csv"
$ws2.Range("B4").Value = "ZB-match_0-16504778364839501.csv"
$ws2.Range("B5").Value = "OB-16504778368799505.csv"
$ws2.Range("B6").Value = "TB-16504778376589866.csv"
$ws2.Range("B7").Value = "OB-16504778371569843.csv"
$ws2.Range("B8").Value = "TB-165047783736795.csv"
$ws2.Range("B9").Value = "TB-1650477838519985.csv"
$ws2.Range("B10").Value = "ZB-match_7-16504778366589482.csv"

# Sheet 3: RS_TO...
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778385479517"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4: TOL_TO...
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778386039479"
$ws4.Range("B2").Value = "MM_stims-16504778385729504.csv"
$ws4.Range("B3").Value = "ZM_stims-1650477838548952.csv"
$ws4.Range("B4").Value = "MM_stims-16504778385879867.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778385729504.csv"
$ws4.Range("B6").Value = "MM_stims-16504778386029851.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477838588949.csv"

# Sheet 5: vSAT_TO...
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778386669817"
$ws5.Range("B2").Value = "vSAT_stims-16504778386509857.csv"
$ws5.Range("B3").Value = "SAT_stims-1650477838618952.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778386349485.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778386079485.csv"
